$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge the two title rows (B2 "SinESC Multi v2.3B" + B3 "Bill of Materials (BOM)")
# into a single title cell, then shift the remaining info rows up by one,
# correcting the component count and its derived consolidation ratio.
$ws.Range("B2").Value = "SinESC Multi v2.3B - Bill of Materials (BOM)"
$ws.Range("B3").Value = "Source: SinESC\Multi Edition\2.3B\SinESC-Multi-2.3B\SinESC-Multi-2.3B.sch"
$ws.Range("B4").Value = "Generated on: 9/3/2020"
$ws.Range("B5").Value = "Tool: Eeschema (5.1.6)-1"
$ws.Range("B6").Value = "Generator: bom_csv_grouped_by_value.py"
$ws.Range("B7").Value = "Total Component Count: 75"
$ws.Range("B8").Value = "Unique Component Count: 23"
$ws.Range("B9").Value = "BOM Consolidation Ratio (Total/Unique): 3.261"
$ws.Range("B10").ClearContents()

# Restore cursor/selection to where the author last clicked.
$ws.Range("B11").Select() | Out-Null
